# Add a new test case row (TC-22) to the Test Cases sheet, matching the
# formatting used by the existing rows, then move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 uses the same style pattern (alignment/border/font) that the new
# row needs, and already has the matching ht=60 row height, so clone its
# formatting down onto the new row before writing values into it.
$fmtSource = $ws.Range("A6:J6")
$newRow = $ws.Range("A24:J24")
$fmtSource.Copy()
$newRow.PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(24, 1).Value = "TC-22"
$ws.Cells.Item(24, 2).Value = "Functionality `nTesting"
$ws.Cells.Item(24, 3).Value = "P1"
$ws.Cells.Item(24, 4).Value = "Verify that Clear Enteries button working fine"
$ws.Cells.Item(24, 5).Value = " application is opened and has been used atleast once"
$ws.Cells.Item(24, 6).Value = "1. Press Clear Enteries Button"
$ws.Cells.Item(24, 7).Value = "All input Fields as well as output fields should get clear"
$ws.Cells.Item(24, 8).Value = "Application worked fine"
$ws.Cells.Item(24, 9).Value = "Pass"
$ws.Cells.Item(24, 10).Value = "Shivank"

$ws.Rows.Item(24).RowHeight = 60

# Move the selection, mirroring where the author's cursor ended up after
# the edit.
$ws.Range("D32").Select() | Out-Null
